# Updates cryptos list: refresh prices, volumes, and reorder some coin rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.092.65"
$ws.Range("E2").Value = "  -1.78%  "

# Row 3
$ws.Range("D3").Value = "1.831.87"

# Row 4
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "'229.60"
$ws.Range("E5").Value = "  -3.88%  "

# Row 6
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").Value = "'0.4628"
$ws.Range("E7").Value = "  -4.18%  "

# Row 8
$ws.Range("D8").Value = "'0.2692"
$ws.Range("E8").Value = "  -6.00%  "

# Row 9
$ws.Range("D9").Value = "'0.06211"
$ws.Range("E9").Value = "  -5.22%  "

# Row 10
$ws.Range("D10").Value = "1.822.16"
$ws.Range("E10").Value = "  -3.53%  "

# Row 11
$ws.Range("D11").Value = "'0.07353"
$ws.Range("E11").Value = "  -1.62%  "

# Row 12
$ws.Range("D12").Value = "'15.98"
$ws.Range("E12").Value = "  -4.27%  "

# Row 13
$ws.Range("D13").Value = "'4.888"
$ws.Range("E13").Value = "  -4.22%  "

# Row 14
$ws.Range("D14").Value = "'82.79"
$ws.Range("E14").Value = "  -6.09%  "

# Row 15
$ws.Range("D15").Value = "'0.6179"
$ws.Range("E15").Value = "  -7.49%  "

# Row 16
$ws.Range("D16").Value = "30.047.62"
$ws.Range("E16").Value = "  -1.84%  "

# Row 17
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("D18").Value = "'227.10"
$ws.Range("E18").Value = "  -2.54%  "

# Row 19
$ws.Range("D19").Value = "'0.000007238"
$ws.Range("E19").Value = "  -4.47%  "

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.11%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.27"
$ws.Range("E21").Value = "  -7.60%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.068.72"
$ws.Range("E22").Value = "  -3.80%  "

# Row 23
$ws.Range("D23").Value = "'4.808"
$ws.Range("E23").Value = "  -8.99%  "

# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.816"
$ws.Range("E24").Value = "  -6.54%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'164.83"
$ws.Range("E25").Value = "  -2.68%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.097"
$ws.Range("E26").Value = "  -2.83%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'17.62"
$ws.Range("E27").Value = "  -6.44%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.839"
$ws.Range("E28").Value = "  -6.27%  "

# Row 29
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.1018"
$ws.Range("E29").Value = "  -0.25%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.369"
$ws.Range("E30").Value = "  -2.07%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.040"
$ws.Range("E31").Value = "  -6.67%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.749"
$ws.Range("E32").Value = "  -6.97%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04784"
$ws.Range("E33").Value = "  -5.60%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.126"
$ws.Range("E34").Value = "  -7.21%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6957"
$ws.Range("E35").Value = "  -7.66%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.677"
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01812"
$ws.Range("E37").Value = "  -3.38%  "

# Row 38
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.609"
$ws.Range("E38").Value = "  -1.52%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.8939"
$ws.Range("E39").Value = "  -2.78%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.905"
$ws.Range("E40").Value = "  -8.09%  "

# Row 41
$ws.Range("D41").Value = "'0.9991"
$ws.Range("E41").Value = "  -0.35%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'103.12"
$ws.Range("E42").Value = "  -3.78%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.450"
$ws.Range("E43").Value = "  -3.47%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3980"
$ws.Range("E44").Value = "  -7.31%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.877"
$ws.Range("E45").Value = "  -7.46%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1183"
$ws.Range("E46").Value = "  -7.15%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'58.93"
$ws.Range("E47").Value = "  -8.33%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.393"
$ws.Range("E48").Value = "  -6.74%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05523"
$ws.Range("E49").Value = "  -2.39%  "

# Row 50
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'32.44"
$ws.Range("E50").Value = "  -4.44%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.349"
$ws.Range("E51").Value = "  -9.87%  "
